$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2 = "2025-10-19T23:56:06.914172"
    3 = "2025-10-19T23:56:06.914172"
    4 = "2025-10-19T23:56:06.915165"
    5 = "2025-10-19T23:56:06.915165"
    6 = "2025-10-19T23:56:06.915165"
    7 = "2025-10-19T23:56:06.915165"
    8 = "2025-10-19T23:56:06.916166"
    9 = "2025-10-19T23:56:06.916166"
    10 = "2025-10-19T23:56:06.916166"
    11 = "2025-10-19T23:56:06.916166"
    12 = "2025-10-19T23:56:06.917163"
    13 = "2025-10-19T23:56:06.917163"
    14 = "2025-10-19T23:56:06.917163"
    15 = "2025-10-19T23:56:06.917163"
    16 = "2025-10-19T23:56:06.918162"
    17 = "2025-10-19T23:56:06.918162"
    18 = "2025-10-19T23:56:06.919163"
    19 = "2025-10-19T23:56:06.919163"
    20 = "2025-10-19T23:56:06.919163"
    21 = "2025-10-19T23:56:06.919163"
    22 = "2025-10-19T23:56:06.919163"
    23 = "2025-10-19T23:56:06.920161"
    24 = "2025-10-19T23:56:06.920161"
    25 = "2025-10-19T23:56:06.920161"
    26 = "2025-10-19T23:56:06.920161"
    27 = "2025-10-19T23:56:06.920161"
    28 = "2025-10-19T23:56:06.921161"
    29 = "2025-10-19T23:56:06.921161"
    30 = "2025-10-19T23:56:06.921161"
    31 = "2025-10-19T23:56:06.921161"
    32 = "2025-10-19T23:56:06.921161"
    33 = "2025-10-19T23:56:06.921161"
    34 = "2025-10-19T23:56:06.922163"
    35 = "2025-10-19T23:56:06.922163"
    36 = "2025-10-19T23:56:06.922163"
    37 = "2025-10-19T23:56:06.922163"
    38 = "2025-10-19T23:56:06.922163"
    39 = "2025-10-19T23:56:06.922163"
    40 = "2025-10-19T23:56:06.923162"
    41 = "2025-10-19T23:56:06.923162"
    42 = "2025-10-19T23:56:06.923162"
    43 = "2025-10-19T23:56:06.923162"
    44 = "2025-10-19T23:56:06.923162"
    45 = "2025-10-19T23:56:06.924160"
    46 = "2025-10-19T23:56:07.022615"
    47 = "2025-10-19T23:56:07.022615"
    48 = "2025-10-19T23:56:07.022615"
    49 = "2025-10-19T23:56:07.023617"
    50 = "2025-10-19T23:56:07.023617"
    51 = "2025-10-19T23:56:07.023617"
    52 = "2025-10-19T23:56:07.024614"
    53 = "2025-10-19T23:56:07.024614"
    54 = "2025-10-19T23:56:07.024614"
    55 = "2025-10-19T23:56:07.024614"
    56 = "2025-10-19T23:56:07.025615"
    57 = "2025-10-19T23:56:07.025615"
    58 = "2025-10-19T23:56:07.025615"
    59 = "2025-10-19T23:56:07.025615"
    60 = "2025-10-19T23:56:07.026614"
    61 = "2025-10-19T23:56:07.026614"
    62 = "2025-10-19T23:56:07.026614"
    63 = "2025-10-19T23:56:07.026614"
    64 = "2025-10-19T23:56:07.026614"
    65 = "2025-10-19T23:56:07.027614"
    66 = "2025-10-19T23:56:07.027614"
    67 = "2025-10-19T23:56:07.027614"
    68 = "2025-10-19T23:56:07.027614"
    69 = "2025-10-19T23:56:07.028617"
    70 = "2025-10-19T23:56:07.028617"
    71 = "2025-10-19T23:56:07.030630"
    72 = "2025-10-19T23:56:07.030630"
    73 = "2025-10-19T23:56:07.031616"
    74 = "2025-10-19T23:56:07.031616"
    75 = "2025-10-19T23:56:07.061402"
    76 = "2025-10-19T23:56:07.061402"
    77 = "2025-10-19T23:56:07.062402"
    78 = "2025-10-19T23:56:07.062402"
    79 = "2025-10-19T23:56:07.063404"
    80 = "2025-10-19T23:56:07.064403"
    81 = "2025-10-19T23:56:07.064403"
    82 = "2025-10-19T23:56:07.064403"
    83 = "2025-10-19T23:56:07.064403"
    84 = "2025-10-19T23:56:07.065408"
    85 = "2025-10-19T23:56:07.065408"
    86 = "2025-10-19T23:56:07.065408"
    87 = "2025-10-19T23:56:07.065408"
    88 = "2025-10-19T23:56:07.065408"
    89 = "2025-10-19T23:56:07.066407"
    90 = "2025-10-19T23:56:07.067402"
    91 = "2025-10-19T23:56:07.067402"
    92 = "2025-10-19T23:56:07.067402"
    93 = "2025-10-19T23:56:07.067402"
    94 = "2025-10-19T23:56:07.067402"
    95 = "2025-10-19T23:56:07.068404"
    96 = "2025-10-19T23:56:07.068404"
    97 = "2025-10-19T23:56:07.068404"
    98 = "2025-10-19T23:56:07.068404"
    99 = "2025-10-19T23:56:07.069404"
    100 = "2025-10-19T23:56:07.069404"
    101 = "2025-10-19T23:56:07.069404"
    102 = "2025-10-19T23:56:07.069404"
    103 = "2025-10-19T23:56:07.090474"
    104 = "2025-10-19T23:56:07.091527"
    105 = "2025-10-19T23:56:07.091527"
    106 = "2025-10-19T23:56:07.092056"
    107 = "2025-10-19T23:56:07.092056"
    108 = "2025-10-19T23:56:07.092056"
    109 = "2025-10-19T23:56:07.092056"
    110 = "2025-10-19T23:56:07.093052"
    111 = "2025-10-19T23:56:07.093052"
    112 = "2025-10-19T23:56:07.093052"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}